$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.619.48"
$ws.Range("E2").Value = "  +3.22%  "
$ws.Range("D3").Value = "1.694.40"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.78"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3945"
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4016"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.518"
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.000"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.22"
$ws.Range("E11").Value = "  -3.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08755"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.220"
$ws.Range("E13").Value = "  +7.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.28"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.206"
$ws.Range("E15").Value = "  +12.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001310"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "1.698.95"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.93"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07065"
$ws.Range("E19").Value = "  +2.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.67"
$ws.Range("E20").Value = "  +3.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.036"
$ws.Range("E21").Value = "  +6.49%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.21"
$ws.Range("E23").Value = "  +3.48%  "
$ws.Range("D24").Value = "24.612.58"
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.128"
$ws.Range("E25").Value = "  +9.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.343"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.76"
$ws.Range("E27").Value = "  +5.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.07"
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "137.01"
$ws.Range("E29").Value = "  +5.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.195"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.456"
$ws.Range("E31").Value = "  +5.25%  "
$ws.Range("D32").Value = "1.886.78"
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.084"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08582"
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.191"
$ws.Range("E35").Value = "  +8.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.61"
$ws.Range("E36").Value = "  +11.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2734"
$ws.Range("E37").Value = "  +3.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.921"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.39"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09123"
$ws.Range("E40").Value = "  +3.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02729"
$ws.Range("E41").Value = "  +8.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.479"
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7638"
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.595"
$ws.Range("E44").Value = "  +8.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7158"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.62"
$ws.Range("E46").Value = "  +3.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.217"
$ws.Range("E47").Value = "  +3.16%  "
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.97"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.329"
$ws.Range("E50").Value = "  +10.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07999"
$ws.Range("E51").Value = "  +2.80%  "
